# Weekly update: insert 5 new price rows at the top of the existing
# "Alcachofa" (artichoke) data block (rows 41-45), shifting the prior
# rows 41-92 down to 46-97. Mirrors a new reporting week (2021-08-05)
# being prepended to the historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 41:92 down by 5 rows (xlShiftDown = -4121) to make
# room for the 5 new rows of data.
$rows = $ws.Rows("41:45")
$rows.Insert(-4121)

$ws.Range("A41").Value = 2
$ws.Range("B41").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C41").Value = "Coquimbo"
$ws.Range("D41").Value = 44413
$ws.Range("E41").Value = 4
$ws.Range("F41").Value = 100112013
$ws.Range("G41").Value = "Alcachofa"
$ws.Range("H41").Value = "Argentina(o)"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 1600
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 9000
$ws.Range("M41").Value = 8500
$ws.Range("N41").Value = "$/caja 50 unidades"
$ws.Range("O41").Value = "Provincia de Limarí"
$ws.Range("P41").Value = 170
$ws.Range("Q41").Value = 50
$ws.Range("R41").Value = "Hortaliza"

$ws.Range("A42").Value = 2
$ws.Range("B42").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C42").Value = "Coquimbo"
$ws.Range("D42").Value = 44413
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = 100112013
$ws.Range("G42").Value = "Alcachofa"
$ws.Range("H42").Value = "Argentina(o)"
$ws.Range("I42").Value = "Segunda"
$ws.Range("J42").Value = 1000
$ws.Range("K42").Value = 6000
$ws.Range("L42").Value = 7000
$ws.Range("M42").Value = 6500
$ws.Range("N42").Value = "$/caja 70 unidades"
$ws.Range("O42").Value = "Provincia de Limarí"
$ws.Range("P42").Value = 93
$ws.Range("Q42").Value = 70
$ws.Range("R42").Value = "Hortaliza"

$ws.Range("A43").Value = 2
$ws.Range("B43").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C43").Value = "Coquimbo"
$ws.Range("D43").Value = 44413
$ws.Range("E43").Value = 4
$ws.Range("F43").Value = 100112013
$ws.Range("G43").Value = "Alcachofa"
$ws.Range("H43").Value = "Española"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 800
$ws.Range("K43").Value = 10000
$ws.Range("L43").Value = 11000
$ws.Range("M43").Value = 10500
$ws.Range("N43").Value = "$/caja 30 unidades"
$ws.Range("O43").Value = "Provincia de Limarí"
$ws.Range("P43").Value = 350
$ws.Range("Q43").Value = 30
$ws.Range("R43").Value = "Hortaliza"

$ws.Range("A44").Value = 2
$ws.Range("B44").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C44").Value = "Coquimbo"
$ws.Range("D44").Value = 44413
$ws.Range("E44").Value = 4
$ws.Range("F44").Value = 100112013
$ws.Range("G44").Value = "Alcachofa"
$ws.Range("H44").Value = "Española"
$ws.Range("I44").Value = "Segunda"
$ws.Range("J44").Value = 800
$ws.Range("K44").Value = 7000
$ws.Range("L44").Value = 8000
$ws.Range("M44").Value = 7500
$ws.Range("N44").Value = "$/caja 40 unidades"
$ws.Range("O44").Value = "Provincia de Limarí"
$ws.Range("P44").Value = 188
$ws.Range("Q44").Value = 40
$ws.Range("R44").Value = "Hortaliza"

$ws.Range("A45").Value = 2
$ws.Range("B45").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C45").Value = "Coquimbo"
$ws.Range("D45").Value = 44413
$ws.Range("E45").Value = 4
$ws.Range("F45").Value = 100112013
$ws.Range("G45").Value = "Alcachofa"
$ws.Range("H45").Value = "Madrigal"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 800
$ws.Range("K45").Value = 9000
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = 9500
$ws.Range("N45").Value = "$/caja 40 unidades"
$ws.Range("O45").Value = "Provincia de Limarí"
$ws.Range("P45").Value = 238
$ws.Range("Q45").Value = 40
$ws.Range("R45").Value = "Hortaliza"

